$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column G values (Направление_подготовки) for rows 2-6
$ws.Range("G2").Value = "09.03.02."
$ws.Range("G3").Value = "09.04.02."
$ws.Range("G4").Value = "09.04.02."
$ws.Range("G5").Value = "09.03.03."
$ws.Range("G6").Value = "09.02.02."

# Update selection / view
$ws.Range("H10").Select()
